$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").Insert(-4121)
Write-Output "done"
